# Auto commit at 2025-11-28 7:55:45.04
# Refresh the monthly metrics figures on the "Metrics" sheet. The "today"
# sheet pulls these same numbers via =Metrics!B.. formulas (and derives
# E/F columns from them), so it recalculates automatically once the
# source values change - no direct writes are needed there.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 363523.51000000007
$wsMetrics.Range("B3").Value  = 320199.15000000002
$wsMetrics.Range("B4").Value  = 111991.78
$wsMetrics.Range("B5").Value  = 14810
$wsMetrics.Range("B6").Value  = 5159769.2600000007
$wsMetrics.Range("B7").Value  = 4362275.830000001
$wsMetrics.Range("B8").Value  = 1518951.6100000003
$wsMetrics.Range("B9").Value  = 201017
$wsMetrics.Range("B10").Value = 33625150.250000015
$wsMetrics.Range("B11").Value = 31637550.990000006
$wsMetrics.Range("B12").Value = 11800673.649999999
$wsMetrics.Range("B13").Value = 1298647

# Cursor/selection moved on both sheets as part of the edit session.
$wsMetrics.Range("D8").Select()

$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("E8").Select()
